# Update weekly price data (Fruta / hortaliza, semanal)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
$data = @{
    2  = @{ D = 44676; J = 120; K = 4000; L = 4500;  M = 4250; P = 71  }
    3  = @{ D = 44657; J = 100; K = 5000; L = 5500;  M = 5250; P = 88  }
    4  = @{ D = 44494; J = 120; K = 5000; L = 6000;  M = 5500; P = 92  }
    5  = @{ D = 44421; J = 100; K = 8000; L = 9000;  M = 8500; P = 142 }
    6  = @{ D = 44362; J = 120; K = 8000; L = 9000;  M = 8500; P = 142 }
    7  = @{ D = 44589; J = 110; K = 5000; L = 6000;  M = 5500; P = 92  }
    9  = @{ D = 44242; J = 160; K = 5000; L = 5500;  M = 5250; P = 88  }
    10 = @{ D = 44400; J = 120; K = 9000; L = 10000; M = 9500; P = 158 }
    11 = @{ D = 44382; J = 160; K = 7000; L = 8000;  M = 7438; P = 124 }
    13 = @{ D = 44669; J = 130; K = 4500; L = 5000;  M = 4750; P = 79  }
    14 = @{ D = 44603; J = 140; K = 5500; L = 6000;  M = 5750; P = 96  }
    15 = @{ D = 44648; J = 120; K = 6500; L = 7000;  M = 6750; P = 112 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals.D   # D - Fecha
    $ws.Cells.Item($row, 10).Value = $vals.J  # J - Volumen
    $ws.Cells.Item($row, 11).Value = $vals.K  # K - Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals.L  # L - Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals.M  # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $vals.P  # P - Precio $/Kg
}
